$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 311, shifting the existing weekly records (rows
# 311-357) down by one. This mirrors the weekly append pattern used in this
# sheet: a brand-new observation is inserted at the top of the data block and
# every later week's data slides down one row (so what was row 357 becomes
# row 358, etc.).
$ws.Rows("311:311").Insert()

# Populate the newly inserted row 311 with the new weekly record. The
# non-varying descriptive columns (Mercado ID, Mercado, Region, Codreg,
# Categoria ID, Categoria, Variedad, Calidad, Unidad de comercializacion,
# Origen, Kg o Unidades, Clasificacion) match every other row in this block.
$ws.Cells.Item(311, 1).Value = 10
$ws.Cells.Item(311, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(311, 3).Value = "La Araucanía"
$ws.Cells.Item(311, 4).Value = 45180
$ws.Cells.Item(311, 5).Value = 9
$ws.Cells.Item(311, 6).Value = 100114007
$ws.Cells.Item(311, 7).Value = "Jengibre"
$ws.Cells.Item(311, 8).Value = "Sin especificar"
$ws.Cells.Item(311, 9).Value = "Primera"
$ws.Cells.Item(311, 10).Value = 180
$ws.Cells.Item(311, 11).Value = 24000
$ws.Cells.Item(311, 12).Value = 24000
$ws.Cells.Item(311, 13).Value = 24000
$ws.Cells.Item(311, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(311, 15).Value = "Perú"
$ws.Cells.Item(311, 16).Value = 1846
$ws.Cells.Item(311, 17).Value = 13
$ws.Cells.Item(311, 18).Value = "Hortaliza"
